# Update feed logs: append two new rows (34, 35) to Sheet1's log table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = "2024-06-15 14:11:45"
$ws.Cells.Item(34, 4).Value = 200
$ws.Cells.Item(34, 5).Value = 8

# Row 35
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = "2024-06-15 14:11:45"
$ws.Cells.Item(35, 4).Value = 200
$ws.Cells.Item(35, 5).Value = 1
